$wb = $excel.ActiveWorkbook

$wsRepay = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N on the "Repayment schedule" sheet,
# pushing the existing "Late" / heading / "Outstanding" columns one slot to the right.
$wsRepay.Columns("N:N").Insert()

# The newly inserted column inherits the width of its left neighbour (column M).
$wsRepay.Columns("N:N").ColumnWidth = $wsRepay.Columns("M:M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab instead of "Transactions",
# and move its selection.
$wsRepay.Select()
$wsRepay.Range("S9").Select()
